$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2892.0625
$ws.Range("J17").Value = 2892.0625
$ws.Range("L17").Value = 8676.1875
$ws.Range("N17").Value = -9012.1875

$ws.Range("H51").Value = 12297.846
$ws.Range("I51").Value = 16250
$ws.Range("J51").Value = 10541.333
$ws.Range("K51").Value = 16250
$ws.Range("L51").Value = 10541.333
$ws.Range("M51").Value = -15766
$ws.Range("N51").Value = -11509.333

$ws.Range("H74").Value = 75735.336
$ws.Range("I74").Value = 75735.336
$ws.Range("K74").Value = 75735.336
$ws.Range("M74").Value = -74799.336

$ws.Range("H77").Value = 75735.336
$ws.Range("I77").Value = 75735.336
$ws.Range("K77").Value = 378676.68
$ws.Range("M77").Value = -373996.68

$ws.Range("H111").Value = 96884.09
$ws.Range("I111").Value = 2251.2
$ws.Range("J111").Value = 175744.83
$ws.Range("K111").Value = 6753.599999999999
$ws.Range("L111").Value = 527234.49
$ws.Range("M111").Value = -3686.599999999999
$ws.Range("N111").Value = -533368.49

$ws.Range("H132").Value = 1802.6271
$ws.Range("I132").Value = 1802.6271
$ws.Range("K132").Value = 5407.8813
$ws.Range("M132").Value = -2877.8813

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2035.0385
$ws.Range("I45").Value = 1315.6666
$ws.Range("J45").Value = 3653.625
$ws.Range("K45").Value = 1315.6666
$ws.Range("L45").Value = 3653.625
$ws.Range("M45").Value = -938.6666
$ws.Range("N45").Value = -4407.625

$ws.Range("H113").Value = 150000
$ws.Range("J113").Value = 150000
$ws.Range("L113").Value = 150000
$ws.Range("N113").Value = -158678

$ws.Range("H122").Value = 2878
$ws.Range("I122").Value = 2266.6875
$ws.Range("J122").Value = 5323.25
$ws.Range("K122").Value = 6800.0625
$ws.Range("L122").Value = 15969.75
$ws.Range("M122").Value = -4350.0625
$ws.Range("N122").Value = -20869.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 5085.3335
$ws.Range("I82").Value = 5085.3335
$ws.Range("K82").Value = 5085.3335
$ws.Range("M82").Value = -4702.3335

$ws.Range("H85").Value = 5085.3335
$ws.Range("I85").Value = 5085.3335
$ws.Range("K85").Value = 5085.3335
$ws.Range("M85").Value = -3759.3335

$ws.Range("H99").Value = 27227.5
$ws.Range("I99").Value = 27227.5
$ws.Range("K99").Value = 27227.5
$ws.Range("M99").Value = -25729.5

$ws.Range("H107").Value = 1001.7778
$ws.Range("I107").Value = 513.8182
$ws.Range("J107").Value = 1768.5714
$ws.Range("K107").Value = 513.8182
$ws.Range("L107").Value = 1768.5714
$ws.Range("M107").Value = 1406.1818
$ws.Range("N107").Value = -5608.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 10000
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H62").Value = 6728.077
$ws.Range("I62").Value = 3660.6667
$ws.Range("J62").Value = 9357.286
$ws.Range("K62").Value = 3660.6667
$ws.Range("L62").Value = 9357.286
$ws.Range("M62").Value = -3036.6667
$ws.Range("N62").Value = -10605.286

$ws.Range("H65").Value = 6728.077
$ws.Range("I65").Value = 3660.6667
$ws.Range("J65").Value = 9357.286
$ws.Range("K65").Value = 18303.3335
$ws.Range("L65").Value = 46786.43
$ws.Range("M65").Value = -15183.3335
$ws.Range("N65").Value = -53026.43

$ws.Range("H107").Value = 2591.647
$ws.Range("I107").Value = 1646
$ws.Range("K107").Value = 1646
$ws.Range("M107").Value = 274

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 980.8333
$ws.Range("J2").Value = 632
$ws.Range("L2").Value = 3792
$ws.Range("N2").Value = -4018

$ws.Range("H12").Value = 66.166664
$ws.Range("I12").Value = 7.3333335
$ws.Range("J12").Value = 125
$ws.Range("K12").Value = 22.0000005
$ws.Range("L12").Value = 375
$ws.Range("M12").Value = 150.9999995
$ws.Range("N12").Value = -721

$ws.Range("H34").Value = 8532.538
$ws.Range("J34").Value = 10963.4
$ws.Range("L34").Value = 32890.2
$ws.Range("N34").Value = -33058.2

$ws.Range("H39").Value = 8279.846
$ws.Range("I39").Value = 2750
$ws.Range("J39").Value = 9285.272000000001
$ws.Range("K39").Value = 8250
$ws.Range("L39").Value = 27855.816
$ws.Range("M39").Value = -7956
$ws.Range("N39").Value = -28443.816

$ws.Range("H55").Value = 9983.333000000001
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 9983.333000000001
$ws.Range("K55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("M55").Value = 29949.999
$ws.Range("N55").Value = -30303.999

$ws.Range("H75").Value = 3540.4
$ws.Range("J75").Value = 4938
$ws.Range("L75").Value = 14814
$ws.Range("N75").Value = -16810

$ws.Range("H78").Value = 3540.4
$ws.Range("J78").Value = 4938
$ws.Range("L78").Value = 44442
$ws.Range("N78").Value = -54426

$ws.Range("H80").Value = 2999.8
$ws.Range("I80").Value = 1666.6666
$ws.Range("J80").Value = 4999.5
$ws.Range("K80").Value = 4999.9998
$ws.Range("L80").Value = 14998.5
$ws.Range("M80").Value = -4063.9998
$ws.Range("N80").Value = -16870.5

$ws.Range("H83").Value = 2999.8
$ws.Range("I83").Value = 1666.6666
$ws.Range("J83").Value = 4999.5
$ws.Range("K83").Value = 14999.9994
$ws.Range("L83").Value = 44995.5
$ws.Range("M83").Value = -10319.9994
$ws.Range("N83").Value = -54355.5

$ws.Range("H87").Value = 772.5
$ws.Range("I87").Value = 772.5
$ws.Range("K87").Value = 2317.5
$ws.Range("M87").Value = -1069.5

$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 30000
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = -29572
$ws.Range("N88").Value = 0

$ws.Range("H90").Value = 772.5
$ws.Range("I90").Value = 772.5
$ws.Range("K90").Value = 6952.5
$ws.Range("M90").Value = -712.5

$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 30000
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = -28518
$ws.Range("N91").Value = 0

$ws.Range("H140").Value = 2302.0967
$ws.Range("I140").Value = 1754.5
$ws.Range("J140").Value = 2383.2222
$ws.Range("K140").Value = 5263.5
$ws.Range("L140").Value = 7149.6666
$ws.Range("M140").Value = -83.5
$ws.Range("N140").Value = -17509.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 496.33334
$ws.Range("I107").Value = 495
$ws.Range("J107").Value = 499
$ws.Range("K107").Value = 495
$ws.Range("L107").Value = 499
$ws.Range("M107").Value = 1425
$ws.Range("N107").Value = -4339

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 170

$ws.Range("H114").Value = 126663.336
$ws.Range("J114").Value = 126663.336
$ws.Range("L114").Value = 126663.336
$ws.Range("N114").Value = -135341.336

$ws.Range("H122").Value = 6988.643
$ws.Range("I122").Value = 2555.25
$ws.Range("J122").Value = 12899.833
$ws.Range("K122").Value = 7665.75
$ws.Range("L122").Value = 38699.499
$ws.Range("M122").Value = -5215.75
$ws.Range("N122").Value = -43599.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6473.095
$ws.Range("J68").Value = 11401.625
$ws.Range("L68").Value = 11401.625
$ws.Range("N68").Value = -12899.625

$ws.Range("H71").Value = 6473.095
$ws.Range("J71").Value = 11401.625
$ws.Range("L71").Value = 57008.125
$ws.Range("N71").Value = -64496.125

$ws.Range("H136").Value = 9749.379000000001
$ws.Range("I136").Value = 6495
$ws.Range("K136").Value = 19485
$ws.Range("M136").Value = -16935

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14997
$ws.Range("I62").Value = 9999.333000000001
$ws.Range("J62").Value = 29990
$ws.Range("K62").Value = 9999.333000000001
$ws.Range("L62").Value = 29990
$ws.Range("M62").Value = -9375.333000000001
$ws.Range("N62").Value = -31238

$ws.Range("H65").Value = 14997
$ws.Range("I65").Value = 9999.333000000001
$ws.Range("J65").Value = 29990
$ws.Range("K65").Value = 49996.665
$ws.Range("L65").Value = 149950
$ws.Range("M65").Value = -46876.665
$ws.Range("N65").Value = -156190

$ws.Range("H126").Value = 3201.2964
$ws.Range("I126").Value = 3201.2964
$ws.Range("K126").Value = 9603.889200000001
$ws.Range("M126").Value = -7133.889200000001

$ws.Range("H136").Value = 7245.0625
$ws.Range("I136").Value = 5763.385
$ws.Range("K136").Value = 17290.155
$ws.Range("M136").Value = -14740.155
